$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(0.2305389221556886,0.5209580838323353,0,0,0,0,0,0,0.02395209580838323,0,0,0,0,0,0.1467065868263473,0,0,0.07784431137724551),
  @(0,0.01694915254237288,0,0,0,0,0,0,0.03954802259887006,0,0,0,0,0,0.7457627118644068,0,0,0.1977401129943503),
  @(0,0,0,0,0,0,0,0,0.1458333333333333,0,0,0,0,0.02083333333333333,0.6875,0,0,0.1458333333333333),
  @(0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0),
  @(0.07602339181286549,0,0.01754385964912281,0,0.07602339181286549,0,0,0,0.2690058479532164,0,0,0,0,0.02339181286549707,0,0.1695906432748538,0.03508771929824561,0.3333333333333333),
  @(0.1069182389937107,0,0.02515723270440252,0,0.06289308176100629,0,0,0,0.1949685534591195,0,0,0,0,0.006289308176100629,0,0.1383647798742138,0.09433962264150944,0.3710691823899371),
  @(0.1188118811881188,0,0.0264026402640264,0,0.03630363036303631,0,0,0,0.1848184818481848,0,0,0,0,0.006600660066006601,0,0.1749174917491749,0.0594059405940594,0.3927392739273927),
  @(0.1344086021505376,0,0.02150537634408602,0.005376344086021506,0.06989247311827956,0,0,0,0.2580645161290323,0,0,0,0,0.005376344086021506,0,0.1451612903225807,0.09139784946236559,0.2688172043010753),
  @(0.1327361563517915,0,0.02361563517915309,0,0.06270358306188925,0,0,0,0.1856677524429967,0,0,0,0,0.01221498371335505,0,0.1872964169381107,0.07573289902280131,0.3200325732899023),
  @(0,0,0,0,0,0.175,0,0,0.07083333333333333,0.2291666666666667,0.5083333333333333,0,0,0,0,0,0,0.01666666666666667),
  @(0,0,0,0,0,0.7768595041322314,0,0,0.1487603305785124,0.01652892561983471,0.02479338842975207,0,0,0,0,0,0,0.03305785123966942),
  @(0,0,0,0,0,0.7631578947368421,0,0,0.2105263157894737,0,0,0,0,0,0,0,0,0.02631578947368421),
  @(0,0,0,0,0,0.75,0,0,0.25,0,0,0,0,0,0,0,0,0),
  @(0,0,0,0,0.005376344086021506,0,0.1989247311827957,0.08064516129032258,0.4247311827956989,0.05376344086021505,0,0.03225806451612903,0,0.03225806451612903,0,0,0,0.1720430107526882),
  @(0,0,0,0,0.01904761904761905,0,0.1333333333333333,0.1238095238095238,0.4142857142857143,0.09523809523809523,0,0.02380952380952381,0.004761904761904762,0.0380952380952381,0,0,0,0.1476190476190476),
  @(0,0,0,0,0.01377410468319559,0,0.1570247933884298,0.09917355371900827,0.440771349862259,0.07713498622589532,0,0.01377410468319559,0,0.05785123966942149,0,0,0,0.140495867768595),
  @(0,0,0,0,0,0,0.1292517006802721,0.08163265306122448,0.4761904761904762,0.1156462585034014,0,0.0272108843537415,0.006802721088435374,0.06802721088435375,0,0,0,0.09523809523809523),
  @(0,0,0,0,0.01867219917012448,0,0.1670124481327801,0.09854771784232365,0.3890041493775934,0.1058091286307054,0,0.01867219917012448,0.002074688796680498,0.0954356846473029,0,0,0,0.1047717842323651)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowVals[$c]
    }
}
